$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the trace-report description text in A1 with the refreshed
# search/completion timestamp (06/15/2023 05:57:14 -> 06/21/2023 08:47:33)
# as part of the new rail car color issue investigation.
$ws.Range("A1").Value = "Description unknown, completed 06/21/2023 08:47:33 EDT, by WPJTOWN1.The search returned: 1 events."
